# Insert a new data row before existing row 75 (Vega Monumental Concepción - Uva),
# shifting all rows from 75 downward down by one row, and populate the
# newly inserted row 75 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75; this shifts rows 75..162 down to 76..163
$ws.Rows.Item(75).Insert()

# Columns A, B, C, E, F, G, H, I, J, L, M keep the same values as the
# (now shifted) Thompson seedless / Primera record that used to be at row 75.
$ws.Cells.Item(75, 1).Value = 11
$ws.Cells.Item(75, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(75, 3).Value = "Bíobío"
$ws.Cells.Item(75, 4).Value = 44740
$ws.Cells.Item(75, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(75, 5).Value = 8
$ws.Cells.Item(75, 6).Value = "Fruta"
$ws.Cells.Item(75, 7).Value = 100109
$ws.Cells.Item(75, 8).Value = "Uva"
$ws.Cells.Item(75, 9).Value = 100109001
$ws.Cells.Item(75, 10).Value = "Uva"
$ws.Cells.Item(75, 11).Value = "Red Globe"
$ws.Cells.Item(75, 12).Value = "Primera"
$ws.Cells.Item(75, 13).Value = 100
$ws.Cells.Item(75, 14).Value = 8000
$ws.Cells.Item(75, 15).Value = 9000
$ws.Cells.Item(75, 16).Value = 8500
$ws.Cells.Item(75, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(75, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(75, 19).Value = 1062
$ws.Cells.Item(75, 20).Value = 8
